# Daily attendance processing - 2025-11-24 10:28:10
# Rotate the "Recorded By" (column G) list of names/emails left by one
# position for every row whose value contains more than one comma-separated
# entry, except rows that mention the admin@admin.com account.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value()

    if ($null -eq $val) { continue }
    if ($val -notmatch ",") { continue }
    if ($val -match "admin@admin\.com") { continue }

    $parts = $val -split ",\s*"
    if ($parts.Count -lt 2) { continue }

    $rotated = $parts[1..($parts.Count - 1)] + $parts[0]
    $cell.Value = [string]::Join(", ", $rotated)
}
